$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "In a dusty old attic, Emily found a magical paintbrush hidden among the cobwebs. With each stroke, the paintbrush brought her drawings to life. She painted colorful landscapes and majestic creatures that danced off the page. Excitedly, Emily shared her creations with her friends and family, filling their lives with wonder and joy. But one day, she realized that her brush could also mend broken hearts and heal old wounds. From that day on, Emily used her gift to spread love and happiness wherever she went, proving that sometimes the most magical things are found in the simplest of places."

$ws.Range("C2").Value = "Em um sótão empoeirado, Emily encontrou um pincel mágico escondido entre as teias de aranha. Com cada traço, o pincel trazia suas desenhos à vida. Ela pintava paisagens coloridas e criaturas majestosas que dançavam fora da página. Animada, Emily compartilhava suas criações com seus amigos e família, enchendo suas vidas de admiração e alegria. Mas um dia, ela percebeu que seu pincel também podia consertar corações partidos e curar velhas feridas. A partir desse dia, Emily usou seu dom para espalhar amor e felicidade por onde passava, provando que às vezes as coisas mais mágicas são encontradas nos lugares mais simples."

$ws.Range("D2").Value = "Básico"

$ws.Range("E2").Value = 609

$ws.Range("G2").Value = "Aprender"
